# feat: add 2022-Q3 data
#
# Before: sheets = [ "总计", "2022-Q2" ]
# After:  sheets = [ "总计", "2022-Q3", "2022-Q2" ]
#   - "2022-Q3" (new, in the slot formerly occupied by "2022-Q2") holds the
#     new quarter's fund-holdings table.
#   - "2022-Q2" (new tab appended after it) keeps the original "2022-Q2"
#     fund-holdings table untouched.
#   - "总计" gains a new row for 2022-Q3 (inserted above the existing
#     2022-Q2 total row).

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q2" sheet so its original data survives on
#    its own tab placed right after it. The duplicate gets appended after
#    the source, matching the sheetId/order the diff expects.
# ---------------------------------------------------------------------------
$q2.Copy($null, $q2)
$q2Copy = $wb.Worksheets.Item($q2.Index + 1)

# The original sheet becomes "2022-Q3" (same tab position / sheetId as
# before), the duplicate becomes the new "2022-Q2".
$q2.Name = "2022-Q3"
$q2Copy.Name = "2022-Q2"
$q3 = $q2

# ---------------------------------------------------------------------------
# 2) Replace the (now-renamed) "2022-Q3" sheet's contents with the new
#    quarter's data. Clear everything below/right of the header row first so
#    stale rows 5-7 (which only existed in the old Q2 table) disappear.
# ---------------------------------------------------------------------------
$q3.Range("A2:H7").Clear()

$q3.Range("B2").Value = "'005994"
$q3.Range("C2").Value = "国投瑞银中证500指数量化增强A"
$q3.Range("D2").Value = "'13.36"
$q3.Range("E2").Value = "'88.67"
$q3.Range("F2").Value = "'1.15"
$q3.Range("G2").Value = "'0.1536"
$q3.Range("H2").Value = 8

$q3.Range("B3").Value = "'007089"
$q3.Range("C3").Value = "国投瑞银中证500指数量化增强C"
$q3.Range("D3").Value = "'4.45"
$q3.Range("E3").Value = "'88.67"
$q3.Range("F3").Value = "'1.15"
$q3.Range("G3").Value = "'0.0512"
$q3.Range("H3").Value = 8

$q3.Range("B4").Value = "'000270"
$q3.Range("C4").Value = "建信灵活配置混合"
$q3.Range("D4").Value = "'2.27"
$q3.Range("E4").Value = "'94.21"
$q3.Range("F4").Value = "'0.93"
$q3.Range("G4").Value = "'0.0211"
$q3.Range("H4").Value = 5

$q3.Range("A2").Value = 0
$q3.Range("A3").Value = 1
$q3.Range("A4").Value = 2

# Re-apply the bold/centered "A2" number style to the refreshed rows (the
# Clear() above also clears formatting), matching the header row's look.
$q3.Range("A2").Copy()
$q3.Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) "总计": add the 2022-Q3 total row above the existing 2022-Q2 total row.
# ---------------------------------------------------------------------------
# Push the current row 2 (2022-Q2 totals) down to row 3, carrying its style.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 6
$total.Range("D3").Value = 0.23

# Write the new 2022-Q3 totals into row 2.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.23
